$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2727
$ws1.Range("F8").Value = 1611
$ws1.Range("F9").Value = 7397
$ws1.Range("F11").Value = 7578
$ws1.Range("F12").Value = 16
$ws1.Range("F15").Value = 6057
$ws1.Range("F16").Value = 3236
$ws1.Range("F17").Value = 3601
$ws1.Range("F19").Value = 3
$ws1.Range("F20").Value = 12
$ws1.Range("F24").Value = 279
$ws1.Range("F25").Value = 276
$ws1.Range("F26").Value = 2090
$ws1.Range("F31").Value = 1059
$ws1.Range("F33").Value = 11
$ws1.Range("G33").Value = 80
$ws1.Range("F37").Value = 8
$ws1.Range("F38").Value = 11
$ws1.Range("F39").Value = 3185
$ws1.Range("F45").Value = 1246
$ws1.Range("F48").Value = 580

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 2727
$ws4.Range("F10").Value = 1611
$ws4.Range("F13").Value = 7397
$ws4.Range("F14").Value = 7578
$ws4.Range("F15").Value = 16
$ws4.Range("F17").Value = 6057
$ws4.Range("F18").Value = 3236
$ws4.Range("F19").Value = 3601
$ws4.Range("F21").Value = 3
$ws4.Range("F26").Value = 279
$ws4.Range("F28").Value = 276
$ws4.Range("F29").Value = 2090
$ws4.Range("F37").Value = 11
$ws4.Range("G37").Value = 80
$ws4.Range("F41").Value = 8
$ws4.Range("F43").Value = 3185
$ws4.Range("F47").Value = 1246
$ws4.Range("F49").Value = 580
